$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in column H, matching the formatting of the other
# header cells (e.g. G1: bold, bordered, centered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the Save values for each data row.
$saveValues = @(0, 0, 1, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
